$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$title1 = $s1.Shapes.Item(1).TextFrame.TextRange
$title1.Text = "."
$title1.Text = "Example numbering MWE"

$s2 = $p.Slides.Item(2)
$title2 = $s2.Shapes.Item(1).TextFrame.TextRange
$title2.Text = "."
$title2.Text = "A second slide"
